$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header label for the "student/gradebook number" column
# -> now references the short code "N_ZACHET"
$ws.Range("C3").Value = "N_ZACHET"

# Match the author's resulting active selection (A4)
$ws.Range("A4").Select()
